# Auto-generated PowerShell Excel COM-interop script
# Updates cryptos list D (Price) and E (Volume 1h) columns per commit diff

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Column D (Price) updates ---
# Values are written with a leading apostrophe to force text interpretation
# (these numeric-looking strings like "5.512" or "28.018.99" must remain text,
# not be auto-converted to numbers/dates by Excel). The cell Style is saved and
# restored afterwards so no extra formatting (e.g. quote-prefix style) is left behind.
$dUpdates = @{
    'D2' = '28.018.99'
    'D3' = '1.911.14'
    'D4' = '1.006'
    'D5' = '315.49'
    'D7' = '0.4821'
    'D8' = '0.3812'
    'D9' = '0.07361'
    'D10' = '0.9344'
    'D11' = '20.81'
    'D12' = '0.07813'
    'D13' = '1.882.53'
    'D14' = '5.512'
    'D15' = '6.637'
    'D16' = '92.09'
    'D18' = '0.000008872'
    'D20' = '28.050.55'
    'D21' = '14.75'
    'D22' = '5.173'
    'D23' = '2.135.99'
    'D25' = '157.19'
    'D26' = '1.918'
    'D28' = '2.137'
    'D29' = '117.08'
    'D30' = '4.974'
    'D31' = '0.08953'
    'D33' = '1.257'
    'D34' = '0.7745'
    'D35' = '4.666'
    'D36' = '2.611'
    'D37' = '0.02049'
    'D38' = '1.109'
    'D39' = '0.5523'
    'D40' = '0.05301'
    'D42' = '7.043'
    'D43' = '0.1529'
    'D44' = '8.514'
    'D45' = '10.73'
    'D46' = '108.68'
    'D47' = '0.4833'
    'D48' = '1.005'
    'D50' = '68.08'
}

foreach ($addr in $dUpdates.Keys) {
    $cell = $ws.Range($addr)
    $origStyle = $cell.Style
    $cell.Value = "'" + $dUpdates[$addr]
    $cell.Style = $origStyle
}

# --- Column E (Volume 1h) updates ---
# These already contain leading/trailing spaces, so Excel keeps them as text naturally.
$eUpdates = @{
    'E2' = '  +2.12%  '
    'E3' = '  +2.58%  '
    'E4' = '  -0.56%  '
    'E5' = '  +1.31%  '
    'E6' = '  -0.60%  '
    'E7' = '  +0.92%  '
    'E8' = '  +0.37%  '
    'E9' = '  +0.53%  '
    'E11' = '  +0.38%  '
    'E12' = '  +0.07%  '
    'E13' = '  +0.79%  '
    'E14' = '  +1.48%  '
    'E15' = '  +1.29%  '
    'E16' = '  +1.93%  '
    'E17' = '  -0.60%  '
    'E18' = '  +0.73%  '
    'E19' = '  -0.52%  '
    'E20' = '  +2.00%  '
    'E21' = '  +0.74%  '
    'E22' = '  +1.11%  '
    'E23' = '  +0.48%  '
    'E24' = '  +2.02%  '
    'E25' = '  +0.89%  '
    'E26' = '  -1.16%  '
    'E27' = '  +0.29%  '
    'E28' = '  +5.90%  '
    'E29' = '  +1.59%  '
    'E30' = '  +0.85%  '
    'E31' = '  +0.83%  '
    'E32' = '  -1.16%  '
    'E33' = '  +3.59%  '
    'E34' = '  +2.20%  '
    'E35' = '  +1.50%  '
    'E36' = '  -4.14%  '
    'E37' = '  +0.15%  '
    'E38' = '  -1.12%  '
    'E39' = '  -0.99%  '
    'E40' = '  +0.63%  '
    'E41' = '  +0.17%  '
    'E42' = '  -0.06%  '
    'E43' = '  +0.31%  '
    'E44' = '  -1.26%  '
    'E45' = '  +0.90%  '
    'E46' = '  +5.51%  '
    'E47' = '  -1.20%  '
    'E48' = '  -0.61%  '
    'E49' = '  -0.25%  '
    'E50' = '  +0.99%  '
    'E51' = '  -0.12%  '
}

foreach ($addr in $eUpdates.Keys) {
    $ws.Range($addr).Value = $eUpdates[$addr]
}
